$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Range("F2").Value = 1.01
$ws.Range("G2").Value = 1.01
$ws.Range("H2").Value = 11.5
$ws.Range("I2").Value = 1000
$ws.Range("J2").Value = 100
$ws.Range("K2").Value = 1000
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = 0
$ws.Range("N2").Value = 0
$ws.Range("O2").Value = 0
$ws.Range("P2").Value = 0
$ws.Range("Q2").Value = 0
$ws.Range("R2").Value = 0
$ws.Range("S2").Value = 0
$ws.Range("T2").Value = 0
$ws.Range("U2").Value = 0
$ws.Range("V2").Value = 1.01
$ws.Range("W2").Value = 100
$ws.Range("X2").Value = 1000
$ws.Range("Y2").Value = 1000
$ws.Range("Z2").Value = 1000
$ws.Range("AA2").Value = 1000
$ws.Range("AB2").Value = 1000
$ws.Range("AC2").Value = 1000
$ws.Range("AD2").Value = 1000
$ws.Range("AE2").Value = 1000
$ws.Range("AF2").Value = 1000
$ws.Range("AG2").Value = 1000
$ws.Range("AH2").Value = 1000
$ws.Range("AI2").Value = 1000
$ws.Range("AJ2").Value = 1000
$ws.Range("AK2").Value = 1000
$ws.Range("AL2").Value = 1000
$ws.Range("AM2").Value = 1000
$ws.Range("AN2").Value = 1.01
$ws.Range("AO2").Value = 1000
# Row 3
$ws.Range("F3").Value = 24
$ws.Range("G3").Value = 32
$ws.Range("H3").Value = 1.1
$ws.Range("I3").Value = 1.11
$ws.Range("J3").Value = 13
$ws.Range("K3").Value = 15.5
$ws.Range("N3").Value = 15
$ws.Range("O3").Value = 1.04
$ws.Range("P3").Value = 6.2
$ws.Range("R3").Value = 3.15
$ws.Range("S3").Value = 1.41
$ws.Range("T3").Value = 1.04
$ws.Range("U3").Value = 1.04
$ws.Range("V3").Value = 8.4
$ws.Range("X3").Value = 230
$ws.Range("Y3").Value = 1000
$ws.Range("Z3").Value = 1000
$ws.Range("AA3").Value = 1000
$ws.Range("AC3").Value = 1000
$ws.Range("AD3").Value = 1000
$ws.Range("AE3").Value = 1000
$ws.Range("AH3").Value = 1000
$ws.Range("AI3").Value = 1000
$ws.Range("AO3").Value = 2.28
# Row 4
$ws.Range("F4").Value = 2.42
$ws.Range("G4").Value = 2.62
$ws.Range("H4").Value = 3.45
$ws.Range("I4").Value = 3.85
$ws.Range("J4").Value = 2.92
$ws.Range("K4").Value = 3.2
$ws.Range("L4").Value = 1.72
$ws.Range("M4").Value = 1.16
$ws.Range("N4").Value = 2.32
$ws.Range("O4").Value = 1.72
$ws.Range("P4").Value = 1.41
$ws.Range("Q4").Value = 3.2
$ws.Range("R4").Value = 1.13
$ws.Range("S4").Value = 7.2
$ws.Range("T4").Value = 2.5
$ws.Range("V4").Value = 1.35
$ws.Range("W4").Value = 1.63
$ws.Range("X4").Value = 7.6
$ws.Range("Y4").Value = 8.6
$ws.Range("Z4").Value = 110
$ws.Range("AA4").Value = 90
$ws.Range("AB4").Value = 6.6
$ws.Range("AC4").Value = 8.199999999999999
$ws.Range("AD4").Value = 32
$ws.Range("AE4").Value = 240
$ws.Range("AF4").Value = 14.5
$ws.Range("AG4").Value = 15
$ws.Range("AJ4").Value = 60
$ws.Range("AL4").Value = 220
$ws.Range("AN4").Value = 600
# Row 5
$ws.Range("F5").Value = 2.06
$ws.Range("G5").Value = 2.12
$ws.Range("H5").Value = 3.3
$ws.Range("I5").Value = 3.5
$ws.Range("J5").Value = 4.3
$ws.Range("K5").Value = 4.5
$ws.Range("L5").Value = 1.27
$ws.Range("N5").Value = 6.8
$ws.Range("O5").Value = 1.14
$ws.Range("P5").Value = 3
$ws.Range("Q5").Value = 1.47
$ws.Range("R5").Value = 1.8
$ws.Range("S5").Value = 2.14
$ws.Range("U5").Value = 2.86
$ws.Range("V5").Value = 1.4
$ws.Range("W5").Value = 1.89
$ws.Range("X5").Value = 32
$ws.Range("Y5").Value = 25
$ws.Range("Z5").Value = 32
$ws.Range("AA5").Value = 65
$ws.Range("AB5").Value = 18.5
$ws.Range("AC5").Value = 11.5
$ws.Range("AD5").Value = 15.5
$ws.Range("AE5").Value = 32
$ws.Range("AF5").Value = 19
$ws.Range("AG5").Value = 11.5
$ws.Range("AI5").Value = 34
$ws.Range("AJ5").Value = 27
$ws.Range("AK5").Value = 18
$ws.Range("AL5").Value = 24
$ws.Range("AM5").Value = 75
$ws.Range("AN5").Value = 8.800000000000001
$ws.Range("AO5").Value = 18
# Row 6
$ws.Range("F6").Value = 2.94
$ws.Range("G6").Value = 3.4
$ws.Range("H6").Value = 2.92
$ws.Range("I6").Value = 3.35
$ws.Range("J6").Value = 2.62
$ws.Range("N6").Value = 2.28
$ws.Range("V6").Value = 1.43
$ws.Range("W6").Value = 1.42
# Row 7
$ws.Range("F7").Value = 1.42
$ws.Range("G7").Value = 1.44
$ws.Range("H7").Value = 8.800000000000001
$ws.Range("I7").Value = 10.5
$ws.Range("J7").Value = 5
$ws.Range("K7").Value = 5.4
$ws.Range("L7").Value = 1.36
$ws.Range("P7").Value = 2.18
$ws.Range("Q7").Value = 1.77
$ws.Range("S7").Value = 2.96
$ws.Range("T7").Value = 2
$ws.Range("U7").Value = 1.84
$ws.Range("V7").Value = 1.11
$ws.Range("W7").Value = 3.25
$ws.Range("X7").Value = 22
$ws.Range("Z7").Value = 85
$ws.Range("AA7").Value = 370
$ws.Range("AB7").Value = 8.6
$ws.Range("AE7").Value = 150
$ws.Range("AF7").Value = 8.4
$ws.Range("AG7").Value = 10.5
$ws.Range("AH7").Value = 28
$ws.Range("AI7").Value = 140
$ws.Range("AL7").Value = 38
$ws.Range("AM7").Value = 170
$ws.Range("AN7").Value = 7
$ws.Range("AO7").Value = 190
# Row 8
$ws.Range("F8").Value = 2.78
$ws.Range("H8").Value = 2.66
$ws.Range("I8").Value = 2.82
$ws.Range("R8").Value = 1.33
$ws.Range("V8").Value = 1.54
$ws.Range("W8").Value = 1.52
$ws.Range("AD8").Value = 14.5
$ws.Range("AH8").Value = 20
$ws.Range("AN8").Value = 600
# Row 9
$ws.Range("F9").Value = 1.51
$ws.Range("G9").Value = 1.55
$ws.Range("H9").Value = 9.4
$ws.Range("I9").Value = 12.5
$ws.Range("J9").Value = 3.75
$ws.Range("K9").Value = 4.4
$ws.Range("L9").Value = 1.63
$ws.Range("M9").Value = 1.12
$ws.Range("N9").Value = 2.48
$ws.Range("O9").Value = 1.59
$ws.Range("P9").Value = 1.51
$ws.Range("Q9").Value = 2.66
$ws.Range("R9").Value = 1.17
$ws.Range("S9").Value = 6
$ws.Range("T9").Value = 2.74
$ws.Range("U9").Value = 1.45
$ws.Range("W9").Value = 2.8
$ws.Range("X9").Value = 14.5
$ws.Range("Y9").Value = 990
$ws.Range("AB9").Value = 5.1
$ws.Range("AC9").Value = 17.5
$ws.Range("AF9").Value = 9
$ws.Range("AG9").Value = 23
$ws.Range("AJ9").Value = 980
$ws.Range("AN9").Value = 600
# Row 10
$ws.Range("F10").Value = 1.66
$ws.Range("H10").Value = 7.2
$ws.Range("I10").Value = 7.4
$ws.Range("J10").Value = 3.75
$ws.Range("K10").Value = 3.8
$ws.Range("N10").Value = 3
$ws.Range("O10").Value = 1.47
$ws.Range("P10").Value = 1.68
$ws.Range("Q10").Value = 2.42
$ws.Range("S10").Value = 4.7
$ws.Range("T10").Value = 2.36
$ws.Range("V10").Value = 1.15
$ws.Range("W10").Value = 2.48
$ws.Range("X10").Value = 9.800000000000001
$ws.Range("Y10").Value = 18
$ws.Range("AA10").Value = 260
$ws.Range("AC10").Value = 8.6
$ws.Range("AD10").Value = 28
$ws.Range("AE10").Value = 140
$ws.Range("AH10").Value = 28
$ws.Range("AM10").Value = 230
$ws.Range("AO10").Value = 240
# Row 11
$ws.Range("H11").Value = 40
$ws.Range("I11").Value = 42
$ws.Range("J11").Value = 12
$ws.Range("K11").Value = 12.5
$ws.Range("L11").Value = 1.23
$ws.Range("N11").Value = 8
$ws.Range("P11").Value = 3.25
$ws.Range("Q11").Value = 1.42
$ws.Range("R11").Value = 1.88
$ws.Range("T11").Value = 2.8
$ws.Range("U11").Value = 1.52
$ws.Range("X11").Value = 48
$ws.Range("Z11").Value = 570
$ws.Range("AB11").Value = 12.5
$ws.Range("AC11").Value = 28
$ws.Range("AF11").Value = 7.4
$ws.Range("AH11").Value = 85
$ws.Range("AJ11").Value = 7.4
$ws.Range("AK11").Value = 16.5
$ws.Range("AM11").Value = 600
$ws.Range("AN11").Value = 2.82
